# Updates TPM-derived values in the Gas6-Mertk sheet (new TPM normalization).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 16.014007
$ws.Cells.Item(2, 8).Value = 48.04202100000001
$ws.Cells.Item(2, 9).Value = 0.09359269140871758
$ws.Cells.Item(2, 10).Value = 0.09359269140871758
$ws.Cells.Item(2, 13).Value = 6.970140000000001
$ws.Cells.Item(2, 14).Value = 20.91042
$ws.Cells.Item(2, 15).Value = 0.06638545261649673
$ws.Cells.Item(2, 16).Value = 0.06638545261649674
$ws.Cells.Item(2, 17).Value = 111.61987075098
$ws.Cells.Item(2, 18).Value = 1004.57883675882
$ws.Cells.Item(2, 19).Value = 0.006213193180763821
$ws.Cells.Item(2, 20).Value = 0.006213193180763822

# Row 3
$ws.Cells.Item(3, 7).Value = 16.014007
$ws.Cells.Item(3, 8).Value = 48.04202100000001
$ws.Cells.Item(3, 9).Value = 0.09359269140871758
$ws.Cells.Item(3, 10).Value = 0.09359269140871758
$ws.Cells.Item(3, 15).Value = 0.03236297878883257
$ws.Cells.Item(3, 16).Value = 0.03236297878883258
$ws.Cells.Item(3, 17).Value = 54.41480576165468
$ws.Cells.Item(3, 18).Value = 489.733251854892
$ws.Cells.Item(3, 19).Value = 0.003028938286850079
$ws.Cells.Item(3, 20).Value = 0.00302893828685008

# Row 4
$ws.Cells.Item(4, 7).Value = 16.014007
$ws.Cells.Item(4, 8).Value = 48.04202100000001
$ws.Cells.Item(4, 9).Value = 0.09359269140871758
$ws.Cells.Item(4, 10).Value = 0.09359269140871758
$ws.Cells.Item(4, 13).Value = 40.232648
$ws.Cells.Item(4, 14).Value = 120.697944
$ws.Cells.Item(4, 15).Value = 0.3831863560043545
$ws.Cells.Item(4, 16).Value = 0.3831863560043545
$ws.Cells.Item(4, 17).Value = 644.2859067005361
$ws.Cells.Item(4, 18).Value = 5798.573160304824
$ws.Cells.Item(4, 19).Value = 0.03586344236954655
$ws.Cells.Item(4, 20).Value = 0.03586344236954655

# Row 5
$ws.Cells.Item(5, 7).Value = 16.014007
$ws.Cells.Item(5, 8).Value = 48.04202100000001
$ws.Cells.Item(5, 9).Value = 0.09359269140871758
$ws.Cells.Item(5, 10).Value = 0.09359269140871758
$ws.Cells.Item(5, 13).Value = 0.6731889999999999
$ws.Cells.Item(5, 14).Value = 2.019567
$ws.Cells.Item(5, 15).Value = 0.006411629674790867
$ws.Cells.Item(5, 16).Value = 0.006411629674790868
$ws.Cells.Item(5, 17).Value = 10.780453358323
$ws.Cells.Item(5, 18).Value = 97.02408022490701
$ws.Cells.Item(5, 19).Value = 0.0006000816775796778
$ws.Cells.Item(5, 20).Value = 0.0006000816775796779

# Row 6
$ws.Cells.Item(6, 7).Value = 16.014007
$ws.Cells.Item(6, 8).Value = 48.04202100000001
$ws.Cells.Item(6, 9).Value = 0.09359269140871758
$ws.Cells.Item(6, 10).Value = 0.09359269140871758
$ws.Cells.Item(6, 13).Value = 53.72106333333333
$ws.Cells.Item(6, 14).Value = 161.16319
$ws.Cells.Item(6, 15).Value = 0.5116535829155252
$ws.Cells.Item(6, 16).Value = 0.5116535829155253
$ws.Cells.Item(6, 17).Value = 860.2894842674434
$ws.Cells.Item(6, 18).Value = 7742.60535840699
$ws.Cells.Item(6, 19).Value = 0.04788703589397744
$ws.Cells.Item(6, 20).Value = 0.04788703589397746

# Row 7
$ws.Cells.Item(7, 9).Value = 0.5957388935007043
$ws.Cells.Item(7, 10).Value = 0.5957388935007044
$ws.Cells.Item(7, 13).Value = 6.970140000000001
$ws.Cells.Item(7, 14).Value = 20.91042
$ws.Cells.Item(7, 15).Value = 0.06638545261649673
$ws.Cells.Item(7, 16).Value = 0.06638545261649674
$ws.Cells.Item(7, 17).Value = 710.4860143779
$ws.Cells.Item(7, 18).Value = 6394.374129401101
$ws.Cells.Item(7, 19).Value = 0.03954839608629519
$ws.Cells.Item(7, 20).Value = 0.03954839608629521

# Row 8
$ws.Cells.Item(8, 9).Value = 0.5957388935007043
$ws.Cells.Item(8, 10).Value = 0.5957388935007044
$ws.Cells.Item(8, 15).Value = 0.03236297878883257
$ws.Cells.Item(8, 16).Value = 0.03236297878883258
$ws.Cells.Item(8, 19).Value = 0.01927988517404588
$ws.Cells.Item(8, 20).Value = 0.01927988517404589

# Row 9
$ws.Cells.Item(9, 9).Value = 0.5957388935007043
$ws.Cells.Item(9, 10).Value = 0.5957388935007044
$ws.Cells.Item(9, 13).Value = 40.232648
$ws.Cells.Item(9, 14).Value = 120.697944
$ws.Cells.Item(9, 15).Value = 0.3831863560043545
$ws.Cells.Item(9, 16).Value = 0.3831863560043545
$ws.Cells.Item(9, 17).Value = 4101.027199652946
$ws.Cells.Item(9, 18).Value = 36909.24479687652
$ws.Cells.Item(9, 19).Value = 0.2282790157306011
$ws.Cells.Item(9, 20).Value = 0.2282790157306012

# Row 10
$ws.Cells.Item(10, 9).Value = 0.5957388935007043
$ws.Cells.Item(10, 10).Value = 0.5957388935007044
$ws.Cells.Item(10, 13).Value = 0.6731889999999999
$ws.Cells.Item(10, 14).Value = 2.019567
$ws.Cells.Item(10, 15).Value = 0.006411629674790867
$ws.Cells.Item(10, 16).Value = 0.006411629674790868
$ws.Cells.Item(10, 17).Value = 68.62005204099832
$ws.Cells.Item(10, 18).Value = 617.580468368985
$ws.Cells.Item(10, 19).Value = 0.003819657167996191
$ws.Cells.Item(10, 20).Value = 0.003819657167996193

# Row 11
$ws.Cells.Item(11, 9).Value = 0.5957388935007043
$ws.Cells.Item(11, 10).Value = 0.5957388935007044
$ws.Cells.Item(11, 13).Value = 53.72106333333333
$ws.Cells.Item(11, 14).Value = 161.16319
$ws.Cells.Item(11, 15).Value = 0.5116535829155252
$ws.Cells.Item(11, 16).Value = 0.5116535829155253
$ws.Cells.Item(11, 17).Value = 5475.93938943016
$ws.Cells.Item(11, 18).Value = 49283.45450487144
$ws.Cells.Item(11, 19).Value = 0.3048119393417659
$ws.Cells.Item(11, 20).Value = 0.304811939341766

# Row 12
$ws.Cells.Item(12, 7).Value = 19.33193133333333
$ws.Cells.Item(12, 8).Value = 57.995794
$ws.Cells.Item(12, 9).Value = 0.1129840572453343
$ws.Cells.Item(12, 10).Value = 0.1129840572453343
$ws.Cells.Item(12, 13).Value = 6.970140000000001
$ws.Cells.Item(12, 14).Value = 20.91042
$ws.Cells.Item(12, 15).Value = 0.06638545261649673
$ws.Cells.Item(12, 16).Value = 0.06638545261649674
$ws.Cells.Item(12, 17).Value = 134.74626786372
$ws.Cells.Item(12, 18).Value = 1212.71641077348
$ws.Cells.Item(12, 19).Value = 0.007500497778679695
$ws.Cells.Item(12, 20).Value = 0.007500497778679697

# Row 13
$ws.Cells.Item(13, 7).Value = 19.33193133333333
$ws.Cells.Item(13, 8).Value = 57.995794
$ws.Cells.Item(13, 9).Value = 0.1129840572453343
$ws.Cells.Item(13, 10).Value = 0.1129840572453343
$ws.Cells.Item(13, 15).Value = 0.03236297878883257
$ws.Cells.Item(13, 16).Value = 0.03236297878883258
$ws.Cells.Item(13, 17).Value = 65.68894896205423
$ws.Cells.Item(13, 18).Value = 591.200540658488
$ws.Cells.Item(13, 19).Value = 0.003656500648107
$ws.Cells.Item(13, 20).Value = 0.003656500648107

# Row 14
$ws.Cells.Item(14, 7).Value = 19.33193133333333
$ws.Cells.Item(14, 8).Value = 57.995794
$ws.Cells.Item(14, 9).Value = 0.1129840572453343
$ws.Cells.Item(14, 10).Value = 0.1129840572453343
$ws.Cells.Item(14, 13).Value = 40.232648
$ws.Cells.Item(14, 14).Value = 120.697944
$ws.Cells.Item(14, 15).Value = 0.3831863560043545
$ws.Cells.Item(14, 16).Value = 0.3831863560043545
$ws.Cells.Item(14, 17).Value = 777.7747884941706
$ws.Cells.Item(14, 18).Value = 6999.973096447535
$ws.Cells.Item(14, 19).Value = 0.04329394918242705
$ws.Cells.Item(14, 20).Value = 0.04329394918242705

# Row 15
$ws.Cells.Item(15, 7).Value = 19.33193133333333
$ws.Cells.Item(15, 8).Value = 57.995794
$ws.Cells.Item(15, 9).Value = 0.1129840572453343
$ws.Cells.Item(15, 10).Value = 0.1129840572453343
$ws.Cells.Item(15, 13).Value = 0.6731889999999999
$ws.Cells.Item(15, 14).Value = 2.019567
$ws.Cells.Item(15, 15).Value = 0.006411629674790867
$ws.Cells.Item(15, 16).Value = 0.006411629674790868
$ws.Cells.Item(15, 17).Value = 13.01404352235533
$ws.Cells.Item(15, 18).Value = 117.126391701198
$ws.Cells.Item(15, 19).Value = 0.0007244119342124556
$ws.Cells.Item(15, 20).Value = 0.0007244119342124557

# Row 16
$ws.Cells.Item(16, 7).Value = 19.33193133333333
$ws.Cells.Item(16, 8).Value = 57.995794
$ws.Cells.Item(16, 9).Value = 0.1129840572453343
$ws.Cells.Item(16, 10).Value = 0.1129840572453343
$ws.Cells.Item(16, 13).Value = 53.72106333333333
$ws.Cells.Item(16, 14).Value = 161.16319
$ws.Cells.Item(16, 15).Value = 0.5116535829155252
$ws.Cells.Item(16, 16).Value = 0.5116535829155253
$ws.Cells.Item(16, 17).Value = 1038.531907513651
$ws.Cells.Item(16, 18).Value = 9346.787167622859
$ws.Cells.Item(16, 19).Value = 0.05780869770190811
$ws.Cells.Item(16, 20).Value = 0.05780869770190812

# Row 17
$ws.Cells.Item(17, 7).Value = 0.6875779999999999
$ws.Cells.Item(17, 8).Value = 2.062734
$ws.Cells.Item(17, 9).Value = 0.004018499278376935
$ws.Cells.Item(17, 10).Value = 0.004018499278376936
$ws.Cells.Item(17, 13).Value = 6.970140000000001
$ws.Cells.Item(17, 14).Value = 20.91042
$ws.Cells.Item(17, 15).Value = 0.06638545261649673
$ws.Cells.Item(17, 16).Value = 0.06638545261649674
$ws.Cells.Item(17, 17).Value = 4.79251492092
$ws.Cells.Item(17, 18).Value = 43.13263428828
$ws.Cells.Item(17, 19).Value = 0.0002667698934341183
$ws.Cells.Item(17, 20).Value = 0.0002667698934341184

# Row 18
$ws.Cells.Item(18, 7).Value = 0.6875779999999999
$ws.Cells.Item(18, 8).Value = 2.062734
$ws.Cells.Item(18, 9).Value = 0.004018499278376935
$ws.Cells.Item(18, 10).Value = 0.004018499278376936
$ws.Cells.Item(18, 15).Value = 0.03236297878883257
$ws.Cells.Item(18, 16).Value = 0.03236297878883258
$ws.Cells.Item(18, 17).Value = 2.336356123485333
$ws.Cells.Item(18, 18).Value = 21.027205111368
$ws.Cells.Item(18, 19).Value = 0.0001300506069090517
$ws.Cells.Item(18, 20).Value = 0.0001300506069090518

# Row 19
$ws.Cells.Item(19, 7).Value = 0.6875779999999999
$ws.Cells.Item(19, 8).Value = 2.062734
$ws.Cells.Item(19, 9).Value = 0.004018499278376935
$ws.Cells.Item(19, 10).Value = 0.004018499278376936
$ws.Cells.Item(19, 13).Value = 40.232648
$ws.Cells.Item(19, 14).Value = 120.697944
$ws.Cells.Item(19, 15).Value = 0.3831863560043545
$ws.Cells.Item(19, 16).Value = 0.3831863560043545
$ws.Cells.Item(19, 17).Value = 27.663083646544
$ws.Cells.Item(19, 18).Value = 248.967752818896
$ws.Cells.Item(19, 19).Value = 0.001539834095087386
$ws.Cells.Item(19, 20).Value = 0.001539834095087386

# Row 20
$ws.Cells.Item(20, 7).Value = 0.6875779999999999
$ws.Cells.Item(20, 8).Value = 2.062734
$ws.Cells.Item(20, 9).Value = 0.004018499278376935
$ws.Cells.Item(20, 10).Value = 0.004018499278376936
$ws.Cells.Item(20, 13).Value = 0.6731889999999999
$ws.Cells.Item(20, 14).Value = 2.019567
$ws.Cells.Item(20, 15).Value = 0.006411629674790867
$ws.Cells.Item(20, 16).Value = 0.006411629674790868
$ws.Cells.Item(20, 17).Value = 0.4628699462419999
$ws.Cells.Item(20, 18).Value = 4.165829516177999
$ws.Cells.Item(20, 19).Value = 0.00002576512922136724
$ws.Cells.Item(20, 20).Value = 0.00002576512922136725

# Row 21
$ws.Cells.Item(21, 7).Value = 0.6875779999999999
$ws.Cells.Item(21, 8).Value = 2.062734
$ws.Cells.Item(21, 9).Value = 0.004018499278376935
$ws.Cells.Item(21, 10).Value = 0.004018499278376936
$ws.Cells.Item(21, 13).Value = 53.72106333333333
$ws.Cells.Item(21, 14).Value = 161.16319
$ws.Cells.Item(21, 15).Value = 0.5116535829155252
$ws.Cells.Item(21, 16).Value = 0.5116535829155253
$ws.Cells.Item(21, 17).Value = 36.93742128460666
$ws.Cells.Item(21, 18).Value = 332.4367915614599
$ws.Cells.Item(21, 19).Value = 0.002056079553725011
$ws.Cells.Item(21, 20).Value = 0.002056079553725012

# Row 22
$ws.Cells.Item(22, 7).Value = 33.136844
$ws.Cells.Item(22, 8).Value = 99.410532
$ws.Cells.Item(22, 9).Value = 0.1936658585668668
$ws.Cells.Item(22, 10).Value = 0.1936658585668668
$ws.Cells.Item(22, 13).Value = 6.970140000000001
$ws.Cells.Item(22, 14).Value = 20.91042
$ws.Cells.Item(22, 15).Value = 0.06638545261649673
$ws.Cells.Item(22, 16).Value = 0.06638545261649674
$ws.Cells.Item(22, 17).Value = 230.9684418381601
$ws.Cells.Item(22, 18).Value = 2078.71597654344
$ws.Cells.Item(22, 19).Value = 0.01285659567732389
$ws.Cells.Item(22, 20).Value = 0.01285659567732389

# Row 23
$ws.Cells.Item(23, 7).Value = 33.136844
$ws.Cells.Item(23, 8).Value = 99.410532
$ws.Cells.Item(23, 9).Value = 0.1936658585668668
$ws.Cells.Item(23, 10).Value = 0.1936658585668668
$ws.Cells.Item(23, 15).Value = 0.03236297878883257
$ws.Cells.Item(23, 16).Value = 0.03236297878883258
$ws.Cells.Item(23, 17).Value = 112.5973611610293
$ws.Cells.Item(23, 18).Value = 1013.376250449264
$ws.Cells.Item(23, 19).Value = 0.006267604072920558
$ws.Cells.Item(23, 20).Value = 0.006267604072920559

# Row 24
$ws.Cells.Item(24, 7).Value = 33.136844
$ws.Cells.Item(24, 8).Value = 99.410532
$ws.Cells.Item(24, 9).Value = 0.1936658585668668
$ws.Cells.Item(24, 10).Value = 0.1936658585668668
$ws.Cells.Item(24, 13).Value = 40.232648
$ws.Cells.Item(24, 14).Value = 120.697944
$ws.Cells.Item(24, 15).Value = 0.3831863560043545
$ws.Cells.Item(24, 16).Value = 0.3831863560043545
$ws.Cells.Item(24, 17).Value = 1333.182980482912
$ws.Cells.Item(24, 18).Value = 11998.64682434621
$ws.Cells.Item(24, 19).Value = 0.07421011462669239
$ws.Cells.Item(24, 20).Value = 0.07421011462669239

# Row 25
$ws.Cells.Item(25, 7).Value = 33.136844
$ws.Cells.Item(25, 8).Value = 99.410532
$ws.Cells.Item(25, 9).Value = 0.1936658585668668
$ws.Cells.Item(25, 10).Value = 0.1936658585668668
$ws.Cells.Item(25, 13).Value = 0.6731889999999999
$ws.Cells.Item(25, 14).Value = 2.019567
$ws.Cells.Item(25, 15).Value = 0.006411629674790867
$ws.Cells.Item(25, 16).Value = 0.006411629674790868
$ws.Cells.Item(25, 17).Value = 22.307358875516
$ws.Cells.Item(25, 18).Value = 200.766229879644
$ws.Cells.Item(25, 19).Value = 0.001241713765781174
$ws.Cells.Item(25, 20).Value = 0.001241713765781174

# Row 26
$ws.Cells.Item(26, 7).Value = 33.136844
$ws.Cells.Item(26, 8).Value = 99.410532
$ws.Cells.Item(26, 9).Value = 0.1936658585668668
$ws.Cells.Item(26, 10).Value = 0.1936658585668668
$ws.Cells.Item(26, 13).Value = 53.72106333333333
$ws.Cells.Item(26, 14).Value = 161.16319
$ws.Cells.Item(26, 15).Value = 0.5116535829155252
$ws.Cells.Item(26, 16).Value = 0.5116535829155253
$ws.Cells.Item(26, 17).Value = 1780.146495190787
$ws.Cells.Item(26, 18).Value = 16021.31845671708
$ws.Cells.Item(26, 19).Value = 0.09908983042414875
$ws.Cells.Item(26, 20).Value = 0.09908983042414876

